# Apply the WTREGEN.xlsx update:
#  - Data sheet: append a new observation row (2023-07-19 / 45126, 531.135)
#  - SeriesInfo sheet: refresh realtime_start/realtime_end, observation_end,
#    and last_updated metadata fields

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: add new row 95, copying the date style used by the rest
#     of column A (numFmt "YYYY-MM-DD HH:MM:SS") via a format-only paste so
#     the existing style table isn't duplicated ---
$lastRow = $dataSheet.Cells.Item(94, 1).Row
$newRow = $lastRow + 1

$dataSheet.Cells.Item($lastRow, 1).Copy() | Out-Null
$dataSheet.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$dataSheet.Cells.Item($newRow, 1).Value = 45126
$dataSheet.Cells.Item($newRow, 2).Value = 531.135

# --- SeriesInfo sheet: update metadata values ---
# realtime_start / realtime_end / observation_end hold plain "YYYY-MM-DD"
# text in the source file (no cell style), so force text entry (otherwise
# Excel's auto-detect would silently turn them into date serials) and then
# drop the temporary text format back off the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $infoSheet.Range("B3") "2023-07-24"
Set-TextValue $infoSheet.Range("B4") "2023-07-24"
Set-TextValue $infoSheet.Range("B7") "2023-07-19"

# last_updated already fails Excel's date/time auto-detect (trailing UTC
# offset), so it round-trips as text with a plain assignment.
$infoSheet.Range("B14").Value = "2023-07-20 15:34:04-05"
